$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Is Active" column (F) currently stores TRUE/FALSE as text strings.
# Push up the actual boolean values instead of their string representation.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
